# Invulformulier.xlsx — add "Functionaliteit" and "Lifts per uur" columns to Blad1,
# with a new dropdown source list on Blad2.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Blad1")
$ws2 = $wb.Worksheets.Item("Blad2")

# ---------------------------------------------------------------------------
# 1) New header "Functionaliteit" in H1 (bold, like the other headers)
# ---------------------------------------------------------------------------
$ws1.Range("H1").Value = "Functionaliteit"
$ws1.Range("H1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 2) New dropdown-source rows on Blad2: Koeling / Lift: Vuilnis / Lift: Anders
# ---------------------------------------------------------------------------
$ws2.Range("A6").Value = "Koeling"
$ws2.Range("A7").Value = "Lift: Vuilnis"
$ws2.Range("A8").Value = "Lift: Anders"

# ---------------------------------------------------------------------------
# 3) New header "Lifts per uur (indien van toepassing)" in I1 (bold)
# ---------------------------------------------------------------------------
$ws1.Range("I1").Value = "Lifts per uur (indien van toepassing)"
$ws1.Range("I1").Font.Bold = $true

# ---------------------------------------------------------------------------
# 4) "Geen" default item, added to Blad2 as the first entry of that list (A5)
# ---------------------------------------------------------------------------
$ws2.Range("A5").Value = "Geen"

# ---------------------------------------------------------------------------
# 5) Fill column H (rows 2-19) with the default "Geen" value
# ---------------------------------------------------------------------------
$ws1.Range("H2:H19").Value = "Geen"

# ---------------------------------------------------------------------------
# 6) Column widths for the two new columns
# ---------------------------------------------------------------------------
$ws1.Columns.Item(8).ColumnWidth = 14.833333333333332
$ws1.Columns.Item(9).ColumnWidth = 13.5

# ---------------------------------------------------------------------------
# 7) Column I (rows 2-19): numeric 0 with a 0.00 number format
# ---------------------------------------------------------------------------
$ws1.Range("I2:I19").NumberFormat = "0.00"
$ws1.Range("I2:I19").Value = 0

# ---------------------------------------------------------------------------
# 8) Data validation (dropdown list) on H2:H19, sourced from Blad2!$A$5:$A$8
# ---------------------------------------------------------------------------
$rangeH = $ws1.Range("H2:H19")
$rangeH.Validation.Add(3, 1, 1, "=Blad2!`$A`$5:`$A`$8")
$rangeH.Validation.IgnoreBlank = $true
$rangeH.Validation.InCellDropdown = $true
$rangeH.Validation.ShowInput = $true
$rangeH.Validation.ShowError = $true

# ---------------------------------------------------------------------------
# 9) Restore cursor/selection state on both sheets
# ---------------------------------------------------------------------------
[void]$ws2.Range("L17").Select()
[void]$ws1.Select()
[void]$ws1.Range("K20").Select()

[void]$wb.Save()
